# Insert a new record row at row 27 (a new daily price observation),
# shifting all existing rows 27-116 down to 28-117.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(27).Insert()

$ws.Range("A27").Value = 4
$ws.Range("B27").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C27").Value = "Los Lagos"
$ws.Range("D27").Value = "2021-10-28"
$ws.Range("E27").Value = 10
$ws.Range("F27").Value = 100112009
$ws.Range("G27").Value = "Acelga"
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 100
$ws.Range("K27").Value = 3000
$ws.Range("L27").Value = 3000
$ws.Range("M27").Value = 3000
$ws.Range("N27").Value = "`$/docena de atados (4 kilos)"
$ws.Range("O27").Value = "Región del Maule"
$ws.Range("P27").Value = 750
$ws.Range("Q27").Value = 4
$ws.Range("R27").Value = "Hortaliza"
